# Shift header values in row 1 so that "variable_trajectory_group" moves
# from column F to column C on every worksheet, pushing the existing
# normalize_group / trajgroup_no_vary_q / uniform_scaling_q headers one
# column to the right (C->D, D->E, E->F).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("C1").Value = "variable_trajectory_group"
    $ws.Range("D1").Value = "normalize_group"
    $ws.Range("E1").Value = "trajgroup_no_vary_q"
    $ws.Range("F1").Value = "uniform_scaling_q"
}
